# Apply the LOTR contacts-sheet update:
#  1. Fix Bilbo Baggins' address from "12 bag end" to "4 bag end".
#  2. Append 10 new vendor rows (15-24) to the contacts table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct Bilbo's address
$ws.Range("D2").Value = "4 bag end"

# 2) New rows of data (Company, First, Last, Address, City, State, Email, Phone)
$newRows = @(
    @("Tom's Farm", "Tom", "Bombadil", "500 Forest Road", "Old Forest", "Eriador", "Tom@tomsfarm.com", "555-222-3333"),
    @("Moria Mining", "Balin", "Son of Fundin", "100 Deep Delve", "Moria", "Rhovanion", "Balin@moriemining.com", "555-222-3333"),
    @("Ent Tree Care", "Treebeard", "Ent", "1 Forest Path", "Fangorn", "Entwood", "Treebeard@entcare.com", "555-444-5555"),
    @("Eagle Transport", "Gwaihir", "Windlord", "1000 Sky Way", "Eyrie", "Misty Mountains", "Gwaihir@eagletransport.com", "555-666-8888"),
    @("Shire Bakery", "Frodo", "Baggins", "4 Bag End", "Westfarthing", "The Shire", "Frodo@shirebakery.com", "555-999-0000"),
    @("Rivendell Music", "Lindir", "Elf", "2 Rivendell Way", "Rivendell", "Eriador", "Lindir@rivendellmusic.com", "555-111-2222"),
    @("Beorn's Honey", "Beorn", "Skin-changer", "300 Honey Lane", "Carrock", "Vale of Anduin", "Beorn@beornshoney.com", "555-333-4444"),
    @("Dale Carpentry", "Bard", "the Bowman", "200 Lake Street", "Dale", "Rhovanion", "Bard@dalecarpentry.com", "555-555-6666"),
    @("Lothlorien Looms", "Galadriel", "Nenya", "1000 Lothlorien Way", "Lothlorien", "Rhovanion", "Galadriel@lothlorienlooms.com", "555-777-8888"),
    @("Bree Inn", "Barliman", "Butterbur", "1 Bree Path", "Bree", "Eriador", "Barliman@breeinn.com", "555-999-1111")
)

$startRow = 15
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
